$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 11181.333
$ws.Range("I9").Value = 12930.429
$ws.Range("J9").Value = 9650.875
$ws.Range("K9").Value = 12930.429
$ws.Range("L9").Value = 9650.875
$ws.Range("M9").Value = -12761.429
$ws.Range("N9").Value = -9988.875
$ws.Range("H98").Value = 797.4643
$ws.Range("I98").Value = 747.4231
$ws.Range("K98").Value = 747.4231
$ws.Range("M98").Value = 750.5769
$ws.Range("H103").Value = 1666.3572
$ws.Range("I103").Value = 1652.4
$ws.Range("K103").Value = 4957.200000000001
$ws.Range("M103").Value = -4371.200000000001
$ws.Range("H122").Value = 797.4643
$ws.Range("I122").Value = 747.4231
$ws.Range("K122").Value = 2242.2693
$ws.Range("M122").Value = 207.7307000000001
$ws.Range("H125").Value = 1516
$ws.Range("I125").Value = 1516
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 13644
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H135").Value = 14498133
$ws.Range("I135").Value = 16672151
$ws.Range("K135").Value = 150049359
$ws.Range("M135").Value = -150046824
$ws.Range("H137").Value = 2068809
$ws.Range("I137").Value = 127200.5
$ws.Range("K137").Value = 381601.5
$ws.Range("M137").Value = -379051.5
$ws.Range("H138").Value = 4079.3274
$ws.Range("J138").Value = 5163.564
$ws.Range("L138").Value = 15490.692
$ws.Range("N138").Value = -25770.692

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 12267.223
$ws.Range("I12").Value = 1300.625
$ws.Range("J12").Value = 100000
$ws.Range("K12").Value = 1300.625
$ws.Range("L12").Value = 100000
$ws.Range("M12").Value = -1127.625
$ws.Range("N12").Value = -100346
$ws.Range("H74").Value = 38475776
$ws.Range("I74").Value = 8569.111000000001
$ws.Range("K74").Value = 8569.111000000001
$ws.Range("M74").Value = -7695.111000000001
$ws.Range("H77").Value = 38475776
$ws.Range("I77").Value = 8569.111000000001
$ws.Range("K77").Value = 42845.55500000001
$ws.Range("M77").Value = -38477.55500000001
$ws.Range("H132").Value = 2856.4138
$ws.Range("I132").Value = 2178.9565
$ws.Range("K132").Value = 6536.869499999999
$ws.Range("M132").Value = -4006.869499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7443.857
$ws.Range("I99").Value = 2214.25
$ws.Range("K99").Value = 2214.25
$ws.Range("M99").Value = -716.25
$ws.Range("H134").Value = 1844.9231
$ws.Range("I134").Value = 1498.6522
$ws.Range("K134").Value = 4495.9566
$ws.Range("M134").Value = -1960.9566

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2625.5
$ws.Range("I10").Value = 1150.8
$ws.Range("K10").Value = 1150.8
$ws.Range("M10").Value = -1011.8
$ws.Range("H16").Value = 3049.7334
$ws.Range("I16").Value = 2854.7
$ws.Range("J16").Value = 3439.8
$ws.Range("K16").Value = 2854.7
$ws.Range("L16").Value = 3439.8
$ws.Range("M16").Value = -2567.7
$ws.Range("N16").Value = -4013.8
$ws.Range("H31").Value = 8371339.5
$ws.Range("I31").Value = 3518412
$ws.Range("K31").Value = 3518412
$ws.Range("M31").Value = -3518117
$ws.Range("H34").Value = 8371339.5
$ws.Range("I34").Value = 3518412
$ws.Range("K34").Value = 3518412
$ws.Range("M34").Value = -3518210
$ws.Range("H86").Value = 460428.22
$ws.Range("J86").Value = 5823.4287
$ws.Range("L86").Value = 5823.4287
$ws.Range("N86").Value = -8069.4287
$ws.Range("H89").Value = 460428.22
$ws.Range("J89").Value = 5823.4287
$ws.Range("L89").Value = 29117.1435
$ws.Range("N89").Value = -40349.14350000001
$ws.Range("H99").Value = 2839
$ws.Range("I99").Value = 2548.75
$ws.Range("K99").Value = 2548.75
$ws.Range("M99").Value = -1050.75
$ws.Range("H113").Value = 3049.7334
$ws.Range("I113").Value = 2854.7
$ws.Range("J113").Value = 3439.8
$ws.Range("K113").Value = 2854.7
$ws.Range("L113").Value = 3439.8
$ws.Range("M113").Value = -684.6999999999998
$ws.Range("N113").Value = -7779.8
$ws.Range("H126").Value = 2839
$ws.Range("I126").Value = 2548.75
$ws.Range("K126").Value = 7646.25
$ws.Range("M126").Value = -5176.25
$ws.Range("H132").Value = 1791.2759
$ws.Range("I132").Value = 1680.2858
$ws.Range("K132").Value = 5040.857400000001
$ws.Range("M132").Value = -2510.857400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 414.22223
$ws.Range("I9").Value = 351.8
$ws.Range("J9").Value = 492.25
$ws.Range("K9").Value = 1055.4
$ws.Range("L9").Value = 1476.75
$ws.Range("M9").Value = -831.4000000000001
$ws.Range("N9").Value = -1924.75
$ws.Range("H68").Value = 3334611
$ws.Range("J68").Value = 3572712
$ws.Range("L68").Value = 10718136
$ws.Range("N68").Value = -10719758
$ws.Range("H69").Value = 650
$ws.Range("I69").Value = 725
$ws.Range("J69").Value = 500
$ws.Range("K69").Value = 2175
$ws.Range("L69").Value = 1500
$ws.Range("M69").Value = -1364
$ws.Range("N69").Value = -3122
$ws.Range("H71").Value = 3334611
$ws.Range("J71").Value = 3572712
$ws.Range("L71").Value = 32154408
$ws.Range("N71").Value = -32162520
$ws.Range("H72").Value = 650
$ws.Range("I72").Value = 725
$ws.Range("J72").Value = 500
$ws.Range("K72").Value = 6525
$ws.Range("L72").Value = 4500
$ws.Range("M72").Value = -2469
$ws.Range("N72").Value = -12612
$ws.Range("H107").Value = 33334208
$ws.Range("J107").Value = 1378
$ws.Range("L107").Value = 4134
$ws.Range("N107").Value = -7974
$ws.Range("H124").Value = 9547.706
$ws.Range("I124").Value = 7476.6665
$ws.Range("K124").Value = 22429.9995
$ws.Range("M124").Value = -17519.9995
$ws.Range("H128").Value = 333124.38
$ws.Range("I128").Value = 333124.38
$ws.Range("K128").Value = 999373.14
$ws.Range("M128").Value = -994393.14
$ws.Range("H137").Value = 6127.5386
$ws.Range("I137").Value = 5565.8
$ws.Range("J137").Value = 8000
$ws.Range("K137").Value = 16697.4
$ws.Range("L137").Value = 24000
$ws.Range("M137").Value = -11597.4
$ws.Range("N137").Value = -34200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 100313.664
$ws.Range("I11").Value = 28974.715
$ws.Range("K11").Value = 28974.715
$ws.Range("M11").Value = -28835.715
$ws.Range("H102").Value = 2145.8
$ws.Range("I102").Value = 1906.8235
$ws.Range("K102").Value = 1906.8235
$ws.Range("M102").Value = -284.8235
$ws.Range("H123").Value = 25842.857
$ws.Range("J123").Value = 25800
$ws.Range("L123").Value = 25800
$ws.Range("N123").Value = -30700
$ws.Range("H132").Value = 19035.451
$ws.Range("I132").Value = 21045.785
$ws.Range("J132").Value = 272.33334
$ws.Range("K132").Value = 63137.355
$ws.Range("L132").Value = 817.0000200000001
$ws.Range("M132").Value = -60607.355
$ws.Range("N132").Value = -5877.00002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2879.1
$ws.Range("I16").Value = 2948.9285
$ws.Range("J16").Value = 2716.1667
$ws.Range("K16").Value = 2948.9285
$ws.Range("L16").Value = 2716.1667
$ws.Range("M16").Value = -2778.9285
$ws.Range("N16").Value = -3056.1667
$ws.Range("H40").Value = 4505
$ws.Range("I40").Value = 3506.2727
$ws.Range("K40").Value = 3506.2727
$ws.Range("M40").Value = -3370.2727
$ws.Range("H46").Value = 1016
$ws.Range("I46").Value = 1016
$ws.Range("K46").Value = 1016
$ws.Range("M46").Value = -828
$ws.Range("H61").Value = 4957.1113
$ws.Range("I61").Value = 4593.5654
$ws.Range("K61").Value = 4593.5654
$ws.Range("M61").Value = -4391.5654
$ws.Range("H113").Value = 4957.1113
$ws.Range("I113").Value = 4593.5654
$ws.Range("K113").Value = 4593.5654
$ws.Range("M113").Value = -2423.5654
$ws.Range("H132").Value = 3035.2886
$ws.Range("I132").Value = 2602.0244
$ws.Range("J132").Value = 4650.1816
$ws.Range("K132").Value = 7806.073199999999
$ws.Range("L132").Value = 13950.5448
$ws.Range("M132").Value = -5276.073199999999
$ws.Range("N132").Value = -19010.5448
$ws.Range("H136").Value = 3778.238
$ws.Range("I136").Value = 3074.611
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 9223.832999999999
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -6673.832999999999
$ws.Range("N136").Value = -29100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H107").Value = 1286.4667
$ws.Range("I107").Value = 1199.8
$ws.Range("K107").Value = 3599.4
$ws.Range("M107").Value = -1679.4
$ws.Range("H109").Value = 40375
$ws.Range("J109").Value = 40375
$ws.Range("L109").Value = 40375
$ws.Range("N109").Value = -43149
$ws.Range("H126").Value = 3321.45
$ws.Range("I126").Value = 2613.2354
$ws.Range("K126").Value = 7839.706200000001
$ws.Range("M126").Value = -5369.706200000001
$ws.Range("H132").Value = 4080.0334
$ws.Range("I132").Value = 4396.385
$ws.Range("J132").Value = 2023.75
$ws.Range("K132").Value = 13189.155
$ws.Range("L132").Value = 6071.25
$ws.Range("M132").Value = -10659.155
$ws.Range("N132").Value = -11131.25
$ws.Range("H136").Value = 1906.05
$ws.Range("I136").Value = 1623.2667
$ws.Range("K136").Value = 4869.800099999999
$ws.Range("M136").Value = -2319.800099999999
